$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the Bangalore branch address (remove stray "xxx" before the
# trailing period) held in cell B4.
$ws.Range("B4").Value = "Mercury Travels Ltd.`n125, Infantry Road, Bangalore 560 001."

# Move the active selection from B4 to C4.
$ws.Range("C4").Select()
